$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '29.364.50'
$ws.Range('E2').Value2 = '  -1.83%  '
$ws.Range('D3').Value2 = '1.852.99'
$ws.Range('E3').Value2 = '  -1.24%  '
$ws.Range('E4').Value2 = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '0.7020'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '238.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  -1.53%  '
$ws.Range('E7').Value2 = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.3074'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -2.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.07378'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '23.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -4.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.08115'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  -3.23%  '
$ws.Range('D12').Value2 = '1.878.32'
$ws.Range('E12').Value2 = '  +0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '0.7266'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  -3.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '5.213'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  -4.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '89.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  -3.28%  '
$ws.Range('D16').Value2 = '29.771.74'
$ws.Range('E16').Value2 = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '5.902'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  -3.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '242.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  -1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.000007738'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  -1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '13.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  -3.61%  '
$ws.Range('B21').Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value2 = '2.178.97'
$ws.Range('E21').Value2 = '  +2.43%  '
$ws.Range('B22').Value2 = 'Dai'
$ws.Range('C22').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '1.005'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +0.68%  '
$ws.Range('E23').Value2 = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '7.638'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  -4.70%  '
$ws.Range('E25').Value2 = '  -4.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '9.041'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value2 = '  -2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '161.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -2.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '18.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  -3.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '1.944'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  -4.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.387'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  -7.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '1.507'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  -1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '4.415'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -3.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '4.068'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  -5.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.05323'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +0.10%  '
$ws.Range('E35').Value2 = '  -3.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.7242'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  -4.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '1.008'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '2.685'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  -0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.01870'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  -4.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '2.721'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  -1.21%  '
$ws.Range('B41').Value2 = 'TrustWalletToken'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.8756'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +2.21%  '
$ws.Range('B42').Value2 = 'TheSandbox'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.4326'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  -3.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '5.940'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  -1.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '69.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -3.64%  '
$ws.Range('E45').Value2 = '  +0.09%  '
$ws.Range('D46').Value2 = '1.032.00'
$ws.Range('E46').Value2 = '  -7.10%  '
$ws.Range('B48').Value2 = 'Aptos'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '7.275'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -4.81%  '
$ws.Range('B49').Value2 = 'RocketPoolETH'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value2 = '2.030.60'
$ws.Range('E49').Value2 = '  +0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '1.760'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  -5.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '9.173'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  -3.27%  '
